$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B3").Value = "Handed back: in sync with en-US"
$wsZh.Range("G2").Value = "2016-02-17 03:34:46"
$wsZh.Range("G3").Value = "2016-02-17 03:34:46"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B3").Value = "Handed back: in sync with en-US"
$wsDe.Range("G2").Value = "2016-02-17 03:35:03"
$wsDe.Range("G3").Value = "2016-02-17 03:35:03"
